# Determine Shipping Price Distance Based.xlsx edit
# - A12 ("Large") becomes "Medium - Large" (price bracket upper bound lowered from 3000 to 450)
# - New row 13 added for "Large" bracket (450 - 3000, price 650)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change label in A12 from "Large" to "Medium - Large" (en dash)
$ws.Range("A12").Value = "Medium " + [char]0x2013 + " Large"

# Update the upper-bound (Distance USL) of the Medium-Large row from 3000 to 450
$ws.Range("C12").Value = 450

# Add new row 13 for the "Large" bracket
$ws.Range("A13").Value = "Large"
$ws.Range("B13").Value = 450
$ws.Range("C13").Value = 3000
$ws.Range("D13").Value = 650
